$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy number-format (date/number) from column D into new columns L:M
# so the newly created cells inherit the correct style, matching the rest of the row.
$ws.Range("D7:D35").Copy()
$ws.Range("L7:M35").PasteSpecial(-4122)
$ws.Range("D38:D77").Copy()
$ws.Range("L38:M77").PasteSpecial(-4122)
$ws.Range("D80:D102").Copy()
$ws.Range("L80:M102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 2: write the refreshed financial data (10 period columns D:M) for every data row.
$arr = New-Object "object[,]" 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D7:M7").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 1022000
$arr[0,1] = 484000
$arr[0,2] = 430000
$arr[0,3] = 374000
$arr[0,4] = 155000
$arr[0,5] = 145000
$arr[0,6] = 350000
$arr[0,7] = 461000
$arr[0,8] = 88000
$arr[0,9] = 251000
$ws.Range("D8:M8").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 26000
$arr[0,1] = 63000
$arr[0,2] = 54000
$arr[0,3] = 39000
$arr[0,4] = 10000
$arr[0,5] = 4000
$arr[0,6] = 8000
$arr[0,7] = 5000
$arr[0,8] = 6000
$arr[0,9] = 31000
$ws.Range("D9:M9").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 996000
$arr[0,1] = 421000
$arr[0,2] = 376000
$arr[0,3] = 335000
$arr[0,4] = 145000
$arr[0,5] = 141000
$arr[0,6] = 342000
$arr[0,7] = 456000
$arr[0,8] = 82000
$arr[0,9] = 220000
$ws.Range("D10:M10").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D11:M11").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 3000
$arr[0,1] = 'NA'
$arr[0,2] = 1000
$arr[0,3] = 2000
$arr[0,4] = 1000
$arr[0,5] = 'NA'
$arr[0,6] = 1000
$arr[0,7] = 1000
$arr[0,8] = 1000
$arr[0,9] = 1000
$ws.Range("D12:M12").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D13:M13").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 18000
$arr[0,1] = 18000
$arr[0,2] = 87000
$arr[0,3] = 17000
$arr[0,4] = 17000
$arr[0,5] = 17000
$arr[0,6] = 15000
$arr[0,7] = 38000
$arr[0,8] = 13000
$arr[0,9] = 12000
$ws.Range("D14:M14").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 226000
$arr[0,1] = 193000
$arr[0,2] = 197000
$arr[0,3] = 161000
$arr[0,4] = 155000
$arr[0,5] = 133000
$arr[0,6] = 141000
$arr[0,7] = 147000
$arr[0,8] = 158000
$arr[0,9] = 150000
$ws.Range("D15:M15").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D16:M16").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 497000
$arr[0,1] = 458000
$arr[0,2] = 504000
$arr[0,3] = 368000
$arr[0,4] = 295000
$arr[0,5] = 172000
$arr[0,6] = 275000
$arr[0,7] = 291000
$arr[0,8] = 310000
$arr[0,9] = 552000
$ws.Range("D17:M17").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 525000
$arr[0,1] = 26000
$arr[0,2] = -74000
$arr[0,3] = 6000
$arr[0,4] = -140000
$arr[0,5] = -27000
$arr[0,6] = 75000
$arr[0,7] = 170000
$arr[0,8] = -222000
$arr[0,9] = -301000
$ws.Range("D18:M18").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D19:M19").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -2000
$arr[0,1] = -2000
$arr[0,2] = 1000
$arr[0,3] = -1000
$arr[0,4] = -1000
$arr[0,5] = 2000
$arr[0,6] = 0
$arr[0,7] = 2000
$arr[0,8] = 1000
$arr[0,9] = 0
$ws.Range("D20:M20").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 749000
$arr[0,1] = 218000
$arr[0,2] = 124000
$arr[0,3] = 173000
$arr[0,4] = 45000
$arr[0,5] = 144000
$arr[0,6] = 246000
$arr[0,7] = 319000
$arr[0,8] = -64000
$arr[0,9] = -151000
$ws.Range("D21:M21").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 40000
$arr[0,1] = 38000
$arr[0,2] = 39000
$arr[0,3] = 46000
$arr[0,4] = 47000
$arr[0,5] = 48000
$arr[0,6] = 46000
$arr[0,7] = 47000
$arr[0,8] = 48000
$arr[0,9] = 49000
$ws.Range("D22:M22").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 483000
$arr[0,1] = -14000
$arr[0,2] = -112000
$arr[0,3] = -41000
$arr[0,4] = -188000
$arr[0,5] = -73000
$arr[0,6] = 29000
$arr[0,7] = 125000
$arr[0,8] = -269000
$arr[0,9] = -350000
$ws.Range("D23:M23").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 130000
$arr[0,1] = -8000
$arr[0,2] = -33000
$arr[0,3] = -15000
$arr[0,4] = -76000
$arr[0,5] = 305000
$arr[0,6] = -298000
$arr[0,7] = 31000
$arr[0,8] = -98000
$arr[0,9] = -132000
$ws.Range("D24:M24").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D25:M25").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 353000
$arr[0,1] = -6000
$arr[0,2] = -79000
$arr[0,3] = -26000
$arr[0,4] = -112000
$arr[0,5] = -378000
$arr[0,6] = 327000
$arr[0,7] = 94000
$arr[0,8] = -171000
$arr[0,9] = -218000
$ws.Range("D26:M26").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 353000
$arr[0,1] = -6000
$arr[0,2] = -83000
$arr[0,3] = -30000
$arr[0,4] = -116000
$arr[0,5] = -381000
$arr[0,6] = 323000
$arr[0,7] = 90000
$arr[0,8] = -174000
$arr[0,9] = -244000
$ws.Range("D27:M27").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D28:M28").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 1000
$arr[0,1] = -1000
$arr[0,2] = -2000
$arr[0,3] = -89000
$arr[0,4] = 74000
$arr[0,5] = 232000
$arr[0,6] = -251000
$arr[0,7] = -2000
$arr[0,8] = -1000
$arr[0,9] = -1000
$ws.Range("D29:M29").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D30:M30").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D31:M31").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 2000
$arr[0,1] = 2000
$arr[0,2] = -1000
$arr[0,3] = 1000
$arr[0,4] = 1000
$arr[0,5] = -2000
$arr[0,6] = 0
$arr[0,7] = -2000
$arr[0,8] = -1000
$arr[0,9] = 0
$ws.Range("D32:M32").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 354000
$arr[0,1] = -7000
$arr[0,2] = -85000
$arr[0,3] = -119000
$arr[0,4] = -42000
$arr[0,5] = -149000
$arr[0,6] = 72000
$arr[0,7] = 88000
$arr[0,8] = -175000
$arr[0,9] = -245000
$ws.Range("D33:M33").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D34:M34").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 354000
$arr[0,1] = -7000
$arr[0,2] = -85000
$arr[0,3] = -119000
$arr[0,4] = -42000
$arr[0,5] = -149000
$arr[0,6] = 72000
$arr[0,7] = 88000
$arr[0,8] = -175000
$arr[0,9] = -245000
$ws.Range("D35:M35").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D38:M38").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D39:M39").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D40:M40").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 3000
$arr[0,1] = 36000
$arr[0,2] = 103000
$arr[0,3] = 681000
$arr[0,4] = 189000
$arr[0,5] = 10000
$arr[0,6] = 8000
$arr[0,7] = 142000
$arr[0,8] = 496000
$arr[0,9] = 623000
$ws.Range("D41:M41").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D42:M42").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 405000
$arr[0,1] = 368000
$arr[0,2] = 323000
$arr[0,3] = 337000
$arr[0,4] = 614000
$arr[0,5] = 268000
$arr[0,6] = 205000
$arr[0,7] = 170000
$arr[0,8] = 168000
$arr[0,9] = 147000
$ws.Range("D43:M43").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 48000
$arr[0,1] = 49000
$arr[0,2] = 40000
$arr[0,3] = 31000
$arr[0,4] = 30000
$arr[0,5] = 42000
$arr[0,6] = 41000
$arr[0,7] = 35000
$arr[0,8] = 36000
$arr[0,9] = 33000
$ws.Range("D44:M44").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 283000
$arr[0,1] = 125000
$arr[0,2] = 162000
$arr[0,3] = 112000
$arr[0,4] = 875000
$arr[0,5] = 328000
$arr[0,6] = 139000
$arr[0,7] = 81000
$arr[0,8] = 54000
$arr[0,9] = 106000
$ws.Range("D45:M45").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 739000
$arr[0,1] = 578000
$arr[0,2] = 628000
$arr[0,3] = 1161000
$arr[0,4] = 1401000
$arr[0,5] = 648000
$arr[0,6] = 393000
$arr[0,7] = 428000
$arr[0,8] = 754000
$arr[0,9] = 909000
$ws.Range("D46:M46").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 167000
$arr[0,1] = 139000
$arr[0,2] = 92000
$arr[0,3] = 85000
$arr[0,4] = 70000
$arr[0,5] = 'NA'
$arr[0,6] = 'NA'
$arr[0,7] = 'NA'
$arr[0,8] = 'NA'
$arr[0,9] = 'NA'
$ws.Range("D47:M47").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 7266000
$arr[0,1] = 7152000
$arr[0,2] = 6974000
$arr[0,3] = 6816000
$arr[0,4] = 6691000
$arr[0,5] = 7384000
$arr[0,6] = 7485000
$arr[0,7] = 7395000
$arr[0,8] = 6474000
$arr[0,9] = 6482000
$ws.Range("D48:M48").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D49:M49").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D50:M50").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D51:M51").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 31000
$arr[0,1] = 49000
$arr[0,2] = 76000
$arr[0,3] = 65000
$arr[0,4] = 45000
$arr[0,5] = 63000
$arr[0,6] = 84000
$arr[0,7] = 87000
$arr[0,8] = 36000
$arr[0,9] = 55000
$ws.Range("D52:M52").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D53:M53").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 8203000
$arr[0,1] = 7918000
$arr[0,2] = 7770000
$arr[0,3] = 8127000
$arr[0,4] = 8207000
$arr[0,5] = 8095000
$arr[0,6] = 7962000
$arr[0,7] = 7910000
$arr[0,8] = 7264000
$arr[0,9] = 7446000
$ws.Range("D54:M54").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D55:M55").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D56:M56").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 514000
$arr[0,1] = 607000
$arr[0,2] = 563000
$arr[0,3] = 522000
$arr[0,4] = 446000
$arr[0,5] = 369000
$arr[0,6] = 348000
$arr[0,7] = 294000
$arr[0,8] = 222000
$arr[0,9] = 185000
$ws.Range("D57:M57").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 'NA'
$arr[0,8] = 'NA'
$arr[0,9] = 125000
$ws.Range("D58:M58").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 201000
$arr[0,1] = 524000
$arr[0,2] = 511000
$arr[0,3] = 403000
$arr[0,4] = 400000
$arr[0,5] = 268000
$arr[0,6] = 271000
$arr[0,7] = 294000
$arr[0,8] = 455000
$arr[0,9] = 314000
$ws.Range("D59:M59").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 715000
$arr[0,1] = 1131000
$arr[0,2] = 1074000
$arr[0,3] = 925000
$arr[0,4] = 846000
$arr[0,5] = 637000
$arr[0,6] = 619000
$arr[0,7] = 588000
$arr[0,8] = 677000
$arr[0,9] = 624000
$ws.Range("D60:M60").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 2485000
$arr[0,1] = 2243000
$arr[0,2] = 2154000
$arr[0,3] = 2576000
$arr[0,4] = 2575000
$arr[0,5] = 2859000
$arr[0,6] = 2601000
$arr[0,7] = 2575000
$arr[0,8] = 2575000
$arr[0,9] = 2574000
$ws.Range("D61:M61").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 702000
$arr[0,1] = 605000
$arr[0,2] = 607000
$arr[0,3] = 620000
$arr[0,4] = 695000
$arr[0,5] = 437000
$arr[0,6] = 438000
$arr[0,7] = 525000
$arr[0,8] = 546000
$arr[0,9] = 614000
$ws.Range("D62:M62").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D63:M63").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D64:M64").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D65:M65").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 3902000
$arr[0,1] = 3979000
$arr[0,2] = 3835000
$arr[0,3] = 4121000
$arr[0,4] = 4080000
$arr[0,5] = 3933000
$arr[0,6] = 3658000
$arr[0,7] = 3688000
$arr[0,8] = 3798000
$arr[0,9] = 3812000
$ws.Range("D66:M66").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D67:M67").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D68:M68").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D69:M69").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 232000
$arr[0,3] = 232000
$arr[0,4] = 232000
$arr[0,5] = 232000
$arr[0,6] = 232000
$arr[0,7] = 232000
$arr[0,8] = 232000
$arr[0,9] = 232000
$ws.Range("D70:M70").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D71:M71").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -3437000
$arr[0,1] = -3791000
$arr[0,2] = -3784000
$arr[0,3] = -3703000
$arr[0,4] = -3588000
$arr[0,5] = -3550000
$arr[0,6] = -3404000
$arr[0,7] = -3480000
$arr[0,8] = -3572000
$arr[0,9] = -3400000
$ws.Range("D72:M72").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D73:M73").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D74:M74").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D75:M75").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 4301000
$arr[0,1] = 3939000
$arr[0,2] = 3703000
$arr[0,3] = 3774000
$arr[0,4] = 3895000
$arr[0,5] = 3930000
$arr[0,6] = 4072000
$arr[0,7] = 3990000
$arr[0,8] = 3234000
$arr[0,9] = 3402000
$ws.Range("D76:M76").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D77:M77").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 43465
$arr[0,1] = 43373
$arr[0,2] = 43281
$arr[0,3] = 43190
$arr[0,4] = 43100
$arr[0,5] = 43008
$arr[0,6] = 42916
$arr[0,7] = 42825
$arr[0,8] = 42735
$arr[0,9] = 42643
$ws.Range("D80:M80").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 354000
$arr[0,1] = -7000
$arr[0,2] = -85000
$arr[0,3] = -119000
$arr[0,4] = -42000
$arr[0,5] = -149000
$arr[0,6] = 72000
$arr[0,7] = 88000
$arr[0,8] = -175000
$arr[0,9] = -245000
$ws.Range("D81:M81").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D82:M82").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 226000
$arr[0,1] = 194000
$arr[0,2] = 197000
$arr[0,3] = 168000
$arr[0,4] = 186000
$arr[0,5] = 169000
$arr[0,6] = 171000
$arr[0,7] = 147000
$arr[0,8] = 157000
$arr[0,9] = 150000
$ws.Range("D83:M83").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D84:M84").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D85:M85").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D86:M86").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D87:M87").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D88:M88").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 231000
$arr[0,1] = 224000
$arr[0,2] = 283000
$arr[0,3] = 145000
$arr[0,4] = 279000
$arr[0,5] = 86000
$arr[0,6] = 120000
$arr[0,7] = 22000
$arr[0,8] = 154000
$arr[0,9] = 29000
$ws.Range("D89:M89").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D90:M90").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -402000
$arr[0,1] = -369000
$arr[0,2] = -356000
$arr[0,3] = -349000
$arr[0,4] = -306000
$arr[0,5] = -313000
$arr[0,6] = -305000
$arr[0,7] = -237000
$arr[0,8] = -276000
$arr[0,9] = -616000
$ws.Range("D91:M91").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D92:M92").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D93:M93").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -493000
$arr[0,1] = -406000
$arr[0,2] = -359000
$arr[0,3] = 362000
$arr[0,4] = 190000
$arr[0,5] = -320000
$arr[0,6] = -275000
$arr[0,7] = -1033000
$arr[0,8] = -150000
$arr[0,9] = -384000
$ws.Range("D94:M94").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = $null
$arr[0,1] = $null
$arr[0,2] = $null
$arr[0,3] = $null
$arr[0,4] = $null
$arr[0,5] = $null
$arr[0,6] = $null
$arr[0,7] = $null
$arr[0,8] = $null
$arr[0,9] = $null
$ws.Range("D95:M95").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D96:M96").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D97:M97").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D98:M98").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D99:M99").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 230000
$arr[0,1] = 116000
$arr[0,2] = -502000
$arr[0,3] = -14000
$arr[0,4] = -290000
$arr[0,5] = 236000
$arr[0,6] = 21000
$arr[0,7] = 657000
$arr[0,8] = -130000
$arr[0,9] = -48000
$ws.Range("D100:M100").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = 0
$arr[0,9] = 0
$ws.Range("D101:M101").Value2 = $arr

$arr = New-Object "object[,]" 1,10
$arr[0,0] = -32000
$arr[0,1] = -66000
$arr[0,2] = -578000
$arr[0,3] = 493000
$arr[0,4] = 179000
$arr[0,5] = 2000
$arr[0,6] = -134000
$arr[0,7] = -354000
$arr[0,8] = -127000
$arr[0,9] = -408000
$ws.Range("D102:M102").Value2 = $arr
